$wb = $excel.ActiveWorkbook

# Rename the "Include from Condition Clinic" sheet to "Include #0"
$includeSheet = $wb.Worksheets.Item("Include from Condition Clinic")
$includeSheet.Name = "Include #0"

# Metadata sheet updates
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B3").Value = "1.0.0"
$metaSheet.Range("B10").Value = "null (http://www.saude.gov.br)"

# Include #0 sheet: translate Portuguese display values to English
$includeSheet.Range("B2").Value = "Active"
$includeSheet.Range("B3").Value = "Recurrence"
$includeSheet.Range("B4").Value = "Relapse"
$includeSheet.Range("B5").Value = "Inactive"
$includeSheet.Range("B6").Value = "Remission"
$includeSheet.Range("B7").Value = "Resolved"
